# Adds perihelion and v0 values to spreadsheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths: new column F (perihelion), widen column G (v0)
# (ColumnWidth values chosen so the saved xlsx <col width=".."> lands on the
#  target 16.13 / 17.0 after Excel's internal pixel-width quantization)
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 15.333333333333332
$ws.Columns.Item(7).ColumnWidth = 16.166666666666664

# ---------------------------------------------------------------------------
# Row 3 headers: perihelion (10^6 km) / v_0 = sqrt(GM/ R_p)
# ---------------------------------------------------------------------------
$ws.Range("F3").Value = "perihelion (10^6 km)"
$ws.Range("F3").HorizontalAlignment = 1
$ws.Range("G3").Value = "v_0 = sqrt(GM/ R_p)"
$ws.Range("G3").HorizontalAlignment = 1

# ---------------------------------------------------------------------------
# Rows 5-12: perihelion values (10^6 km) in column F
# ---------------------------------------------------------------------------
$perihelion = @{
    5  = 46.0
    6  = 107.5
    7  = 147.1
    8  = 206.7
    9  = 740.6
    10 = 1357.6
    11 = 2732.7
    12 = 4471.1
}
foreach ($row in $perihelion.Keys) {
    $cell = $ws.Range("F$row")
    $cell.Value = $perihelion[$row]
    $cell.HorizontalAlignment = 1
}

# ---------------------------------------------------------------------------
# Rows 5-12: v_0 = sqrt(GM/R_p) formula in column G
# ---------------------------------------------------------------------------
foreach ($row in 5..12) {
    $cell = $ws.Range("G$row")
    $cell.Formula = "=SQRT(((6.67*10^(-11)) * (1.989*10^30))/(F$row*10^9))"
    $cell.HorizontalAlignment = 1
}

# ---------------------------------------------------------------------------
# Row 17 headers (scaled table): perihelion / v inital
# ---------------------------------------------------------------------------
$ws.Range("F17").Value = "perihelion"
$ws.Range("F17").HorizontalAlignment = 1
$ws.Range("G17").Value = "v inital"
$ws.Range("G17").HorizontalAlignment = 1

# ---------------------------------------------------------------------------
# Rows 19-26 (scaled table): perihelion *10 in column F, v0 scaled in column G
# ---------------------------------------------------------------------------
$srcRow = @{
    19 = 5
    20 = 6
    21 = 7
    22 = 8
    23 = 9
    24 = 10
    25 = 11
    26 = 12
}
foreach ($row in $srcRow.Keys) {
    $src = $srcRow[$row]

    $fcell = $ws.Range("F$row")
    $fcell.Formula = "=F$src * 10"
    $fcell.HorizontalAlignment = 1

    $gcell = $ws.Range("G$row")
    $gcell.Formula = "=G$src*10^-3"
    $gcell.HorizontalAlignment = 1
}
